$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.566.65'
$ws.Range('E2').Value = '  +4.27%  '
$ws.Range('D3').Value = '2.336.66'
$ws.Range('E3').Value = '  +2.30%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '547.84'
$ws.Range('E5').Value = '  +2.42%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '131.78'
$ws.Range('E6').Value = '  +0.42%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.581'
$ws.Range('E8').Value = '  -0.98%  '
$ws.Range('D9').Value = '2.334.38'
$ws.Range('E9').Value = '  +2.19%  '
$ws.Range('E10').Value = '  +1.52%  '
$ws.Range('E11').Value = '  +0.87%  '
$ws.Range('E13').Value = '  +1.58%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.87'
$ws.Range('E14').Value = '  +1.58%  '
$ws.Range('D15').Value = '2.754.50'
$ws.Range('E15').Value = '  +2.34%  '
$ws.Range('D16').Value = '60.557.62'
$ws.Range('E16').Value = '  +4.37%  '
$ws.Range('E17').Value = '  +1.31%  '
$ws.Range('D18').Value = '2.336.68'
$ws.Range('E18').Value = '  +1.82%  '
$ws.Range('E19').Value = '  +1.35%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '316.10'
$ws.Range('E21').Value = '  +0.94%  '
$ws.Range('E22').Value = '  +3.60%  '
$ws.Range('E23').Value = '  -0.37%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '64.20'
$ws.Range('E24').Value = '  +1.78%  '
$ws.Range('E25').Value = '  +1.54%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  -0.28%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.92'
$ws.Range('E27').Value = '  -0.56%  '
$ws.Range('E28').Value = '  +6.86%  '
$ws.Range('E29').Value = '  +12.50%  '
$ws.Range('E30').Value = '  +1.46%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.74'
$ws.Range('E31').Value = '  +2.16%  '
$ws.Range('D32').Value = '0.0₃0737'
$ws.Range('E32').Value = '  +2.26%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.97'
$ws.Range('E33').Value = '  +3.56%  '
$ws.Range('E34').Value = '  +12.14%  '
$ws.Range('E35').Value = '  +0.20%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '17.95'
$ws.Range('E36').Value = '  +0.77%  '
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('E38').Value = '  +0.00%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.10'
$ws.Range('E39').Value = '  +5.20%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '329.45'
$ws.Range('E40').Value = '  +14.15%  '
$ws.Range('E41').Value = '  +3.64%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '38.09'
$ws.Range('E42').Value = '  -0.60%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '140.09'
$ws.Range('E43').Value = '  +0.22%  '
$ws.Range('E44').Value = '  +1.73%  '
$ws.Range('E45').Value = '  -0.49%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '19.36'
$ws.Range('E46').Value = '  +7.19%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0498'
$ws.Range('E47').Value = '  +0.76%  '
$ws.Range('E48').Value = '  +2.00%  '
$ws.Range('D49').Value = '0.0₆0220'
$ws.Range('E49').Value = '  +20.69%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0214'
$ws.Range('E50').Value = '  +2.08%  '
$ws.Range('E51').Value = '  +0.75%  '
